# Generate Report for Handback
#
# Refresh the handback-status report with the latest xliff generate /
# handback timestamps for the b64c069e-... file (row 2 of each table):
#   - Overview!G2            "Latest HO Xliff Generate Date"
#   - zh-cn!H2 / zh-cn!K2    "Correspond Handoff Datetime" / "Correspond Handback DateTime"
#   - de-de!H2 / de-de!K2    "Correspond Handoff Datetime" / "Correspond Handback DateTime"

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-18 00:45:45"

$wsZhCn.Range("H2").Value = "2016-08-18 00:45:40"
$wsZhCn.Range("K2").Value = "2016-08-18 00:45:57"

$wsDeDe.Range("H2").Value = "2016-08-18 00:45:45"
$wsDeDe.Range("K2").Value = "2016-08-18 00:46:10"
